{"js": "// Locate the CALENDARIO DE LA PROPUESTA table: it is the table with a\n// header row (\"ACTIVIDAD\"/\"INICIO\"/\"TERMINO\"/\"OBSERVACI\u00d3N\") followed by a\n// single templated data row. We find it by scanning the document's tables\n// for one whose first row starts with \"ACTIVIDAD\".\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < tables.items.length; i++) {\n  tables.items[i].load(\"values\");\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < tables.items.length; i++) {\n  const vals = tables.items[i].values;\n  if (vals && vals.length > 0 && vals[0] && vals[0][0] === \"ACTIVIDAD\") {\n    target = tables.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"CALENDARIO table not found\");\n}\n\ntarget.rows.load(\"items\");\nawait context.sync();\n\nconst headerRow = target.rows.items[0];\nheaderRow.cells.load(\"items\");\nawait context.sync();\n\n// --- 1) Resize the four columns (dxa -> points, 20 dxa = 1 pt) ---------\n// Setting TableCell.columnWidth resizes the whole column (every row's cell\n// plus <w:gridCol>), matching the OOXML diff where both the grid and each\n// row's tcW changed together.\nconst newWidthsDxa = [2935, 1171, 1276, 3544];\nfor (let i = 0; i < headerRow.cells.items.length; i++) {\n  headerRow.cells.items[i].columnWidth = newWidthsDxa[i] / 20;\n}\nawait context.sync();\n\n// --- 2) Fix the Jinja template text in the data row ---------------------\nconst dataRow = target.rows.items[1];\ndataRow.cells.load(\"items\");\nawait context.sync();\n\n// Cell 0 (ACTIVIDAD): first paragraph holds the mistyped \"{%tr for i in\n// calendario%}\" loop-open tag; the second paragraph \"{{ i.actividad }}\"\n// stays as-is.\nconst cell0Body = dataRow.cells.items[0].body;\ncell0Body.paragraphs.load(\"items\");\nawait context.sync();\ncell0Body.paragraphs.items[0].insertText(\"{% for i in calendario %}\", \"Replace\");\n\n// Cell 3 (OBSERVACI\u00d3N): single paragraph \"{{ i.obs }} {%entr%}\" where the\n// broken \"{%entr%}\" closing tag becomes \"{% endfor %}\".\nconst cell3Body = dataRow.cells.items[3].body;\ncell3Body.paragraphs.load(\"items\");\nawait context.sync();\ncell3Body.paragraphs.items[0].insertText(\"{{ i.obs }} {% endfor %}\", \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"CALENDARIO DE LA PROPUESTA\" table: the one whose first\n# header cell reads \"ACTIVIDAD\" (cell text carries a trailing cell-mark,\n# hence the wildcard match).\n$target = $null\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables.Item($i)\n    if ($t.Cell(1, 1).Range.Text -like \"ACTIVIDAD*\") {\n        $target = $t\n    }\n}\n\n# --- 1) Resize the four columns (dxa -> points, 20 dxa = 1 pt) ----------\n# Setting Column.Width resizes every cell in that column (both the header\n# and data rows) plus the shared <w:gridCol>, matching the OOXML diff.\n$target.Columns.Item(1).Width = 146.75\n$target.Columns.Item(2).Width = 58.55\n$target.Columns.Item(3).Width = 63.8\n$target.Columns.Item(4).Width = 177.2\n\n# --- 2) Fix the Jinja template text in the data row ----------------------\n# \"{%tr for i in calendario%}\" (mistyped loop-open tag) -> \"{% for i in calendario %}\"\n$rng1 = $d.Content\n$rng1.Find.Execute(\"{%tr for i in calendario%}\", $false, $false, $false, $false, $false, $true, 1, $false, \"{% for i in calendario %}\", 2)\n\n# \"{%entr%}\" (broken loop-close tag) -> \"{% endfor %}\"\n$rng2 = $d.Content\n$rng2.Find.Execute(\"{%entr%}\", $false, $false, $false, $false, $false, $true, 1, $false, \"{% endfor %}\", 2)\n"}
